$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("cor_z")
$ws.Range("B3").Value = -29.24564751498815
$ws.Range("B4").Value = -27.64564751498822
$ws.Range("B5").Value = -26.04564751498592
$ws.Range("B6").Value = -24.44564751498288
$ws.Range("B7").Value = -22.84564751498458
$ws.Range("B8").Value = -21.2456475149837
$ws.Range("B9").Value = -19.64564751498684
$ws.Range("B10").Value = -18.04564751499014
$ws.Range("B11").Value = -16.44564751498885
$ws.Range("B12").Value = -14.84564751498526
$ws.Range("B13").Value = -13.24564751498398
$ws.Range("B14").Value = -11.64564751498859
$ws.Range("B15").Value = -10.04564751498664
$ws.Range("B16").Value = -8.445647514986666
$ws.Range("B17").Value = -6.845647514987343
$ws.Range("B18").Value = -5.245647514989332
$ws.Range("B19").Value = -3.645647514986078
$ws.Range("B20").Value = -2.04564751499462
$ws.Range("B21").Value = -0.4456475150033954

$ws = $wb.Worksheets.Item("mome_y")
$ws.Range("B3").Value = 4940.502055675838
$ws.Range("B4").Value = 4371.589105376085
$ws.Range("B5").Value = 3834.676155076309
$ws.Range("B6").Value = 3329.763204776648
$ws.Range("B7").Value = 2856.850254476994
$ws.Range("B8").Value = 2415.937304177264
$ws.Range("B9").Value = 2007.024353877646
$ws.Range("B10").Value = 1630.111403577791
$ws.Range("B11").Value = 1285.19845327804
$ws.Range("B12").Value = 972.285502978237
$ws.Range("B13").Value = 691.3725526786297
$ws.Range("B14").Value = 442.4596023788779
$ws.Range("B15").Value = 225.5466520790862
$ws.Range("B16").Value = 40.63370177941211
$ws.Range("B17").Value = -112.2792485203804
$ws.Range("B18").Value = -233.1921988200816
$ws.Range("B19").Value = -322.1051491199537
$ws.Range("B20").Value = -379.0180994195247
$ws.Range("B21").Value = -403.9310497197385

$ws = $wb.Worksheets.Item("drz")
$ws.Range("B2").Value = 0.000000000000000001129851574057469
$ws.Range("B3").Value = 0.02346900992551523
$ws.Range("B4").Value = 0.08878606467591761
$ws.Range("B5").Value = 0.1911322421300734
$ws.Range("B6").Value = 0.3259596730718583
$ws.Range("B7").Value = 0.4889915411901583
$ws.Range("B8").Value = 0.6762220830788691
$ws.Range("B9").Value = 0.8839165882368952
$ws.Range("B10").Value = 1.108611399068152
$ws.Range("B11").Value = 1.347113910881561
$ws.Range("B12").Value = 1.596502571891056
$ws.Range("B13").Value = 1.854126883215578
$ws.Range("B14").Value = 2.117607398879081
$ws.Range("B15").Value = 2.384835725810527
$ws.Range("B16").Value = 2.653974523843885
$ws.Range("B17").Value = 2.923457505718136
$ws.Range("B18").Value = 3.191989437077271
$ws.Range("B19").Value = 3.458546136470289
$ws.Range("B20").Value = 3.722374475351199
$ws.Range("B21").Value = 3.982992378079019
